$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.022.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.419.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.148"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.354"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.848.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.965.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.392.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "329.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0778"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("E31").Value = "  -4.14%  "
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "321.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0972"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -3.49%  "
$ws.Range("E51").Value = "  -1.09%  "
